# Update workbook according to the authoritative diff:
#  - "Förändrad" (column C) timestamp bumped from 2023-11-17 (45247) to 2023-11-18 (45248)
#    for every existing data row (2-28).
#  - Rows 2 and 3 swap their data (row 2 now holds "A 57664-2023", row 3 holds "A 57619-2023").
#  - Three brand-new rows (29, 30, 31) are appended with fresh case data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "Förändrad" column (C) for all existing data rows (2-28) to 45248.
$ws.Range("C2:C28").Value2 = 45248

# 2. Swap the contents of rows 2 and 3 (columns A through R), since the new
#    case "A 57664-2023" now sorts above "A 57619-2023".
$row2 = $ws.Range("A2:R2").Value2
$row3 = $ws.Range("A3:R3").Value2
$ws.Range("A2:R2").Value2 = $row3
$ws.Range("A3:R3").Value2 = $row2

# ...and swap the link formulas (S,T,V,W,X,Y) that reference the case id.
$f2 = $ws.Range("S2:Y2").Formula
$f3 = $ws.Range("S3:Y3").Formula
$ws.Range("S2:Y2").Formula = $f3
$ws.Range("S3:Y3").Formula = $f2

# Row 28 picks up an explicit row height once it is no longer the last row.
$ws.Rows.Item(28).RowHeight = 15

# 3. Append new rows 29-31 with the newly reported cases.

# Row 29: A 57804-2023
$ws.Rows.Item(29).RowHeight = 15
$ws.Cells.Item(29, 1).Value2 = "A 57804-2023"
$ws.Cells.Item(29, 2).Value2 = 45243
$ws.Cells.Item(29, 3).Value2 = 45248
$ws.Cells.Item(29, 4).Value2 = "OKÄNT"
$ws.Cells.Item(29, 5).Value2 = "OKÄNT"
$ws.Cells.Item(29, 6).Value2 = "SCA"
$ws.Cells.Item(29, 7).Value2 = 7.3
$ws.Cells.Item(29, 8).Value2 = 0
$ws.Cells.Item(29, 9).Value2 = 0
$ws.Cells.Item(29, 10).Value2 = 0
$ws.Cells.Item(29, 11).Value2 = 0
$ws.Cells.Item(29, 12).Value2 = 0
$ws.Cells.Item(29, 13).Value2 = 0
$ws.Cells.Item(29, 14).Value2 = 0
$ws.Cells.Item(29, 15).Value2 = 0
$ws.Cells.Item(29, 16).Value2 = 0
$ws.Cells.Item(29, 17).Value2 = 0
$ws.Cells.Item(29, 18).WrapText = $true
$ws.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(29, 3).NumberFormat = "YYYY-MM-DD"

# Row 30: A 57805-2023
$ws.Rows.Item(30).RowHeight = 15
$ws.Cells.Item(30, 1).Value2 = "A 57805-2023"
$ws.Cells.Item(30, 2).Value2 = 45243
$ws.Cells.Item(30, 3).Value2 = 45248
$ws.Cells.Item(30, 4).Value2 = "OKÄNT"
$ws.Cells.Item(30, 5).Value2 = "OKÄNT"
$ws.Cells.Item(30, 6).Value2 = "SCA"
$ws.Cells.Item(30, 7).Value2 = 5.6
$ws.Cells.Item(30, 8).Value2 = 0
$ws.Cells.Item(30, 9).Value2 = 0
$ws.Cells.Item(30, 10).Value2 = 0
$ws.Cells.Item(30, 11).Value2 = 0
$ws.Cells.Item(30, 12).Value2 = 0
$ws.Cells.Item(30, 13).Value2 = 0
$ws.Cells.Item(30, 14).Value2 = 0
$ws.Cells.Item(30, 15).Value2 = 0
$ws.Cells.Item(30, 16).Value2 = 0
$ws.Cells.Item(30, 17).Value2 = 0
$ws.Cells.Item(30, 18).WrapText = $true
$ws.Cells.Item(30, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 3).NumberFormat = "YYYY-MM-DD"

# Row 31: A 57807-2023 (no Markägare on record for this one)
$ws.Cells.Item(31, 1).Value2 = "A 57807-2023"
$ws.Cells.Item(31, 2).Value2 = 45244
$ws.Cells.Item(31, 3).Value2 = 45248
$ws.Cells.Item(31, 4).Value2 = "OKÄNT"
$ws.Cells.Item(31, 5).Value2 = "OKÄNT"
$ws.Cells.Item(31, 7).Value2 = 1.4
$ws.Cells.Item(31, 8).Value2 = 0
$ws.Cells.Item(31, 9).Value2 = 0
$ws.Cells.Item(31, 10).Value2 = 0
$ws.Cells.Item(31, 11).Value2 = 0
$ws.Cells.Item(31, 12).Value2 = 0
$ws.Cells.Item(31, 13).Value2 = 0
$ws.Cells.Item(31, 14).Value2 = 0
$ws.Cells.Item(31, 15).Value2 = 0
$ws.Cells.Item(31, 16).Value2 = 0
$ws.Cells.Item(31, 17).Value2 = 0
$ws.Cells.Item(31, 18).WrapText = $true
$ws.Cells.Item(31, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(31, 3).NumberFormat = "YYYY-MM-DD"
